{"js": "// Replace the date header and all the two-digit multiplication problems.\n// Every source string in this document appears exactly once, so a plain\n// body.search() + Range.insertText(\"Replace\") round-trip is safe for each.\nconst replacements = [\n  [\"2024-03-07 Thursday\", \"2024-03-08 Friday\"],\n  [\"40\u00d748=\", \"13\u00d777=\"],\n  [\"74\u00d798=\", \"76\u00d732=\"],\n  [\"94\u00d788=\", \"76\u00d783=\"],\n  [\"85\u00d711=\", \"32\u00d720=\"],\n  [\"30\u00d762=\", \"71\u00d792=\"],\n  [\"57\u00d785=\", \"83\u00d727=\"],\n  [\"76\u00d748=\", \"28\u00d740=\"],\n  [\"17\u00d776=\", \"87\u00d748=\"],\n  [\"53\u00d748=\", \"28\u00d781=\"],\n  [\"45\u00d789=\", \"29\u00d716=\"],\n  [\"97\u00d748=\", \"92\u00d763=\"],\n  [\"88\u00d715=\", \"66\u00d775=\"],\n  [\"14\u00d741=\", \"84\u00d727=\"],\n  [\"67\u00d722=\", \"53\u00d732=\"],\n  [\"80\u00d793=\", \"18\u00d799=\"],\n  [\"42\u00d734=\", \"38\u00d786=\"],\n  [\"62\u00d777=\", \"61\u00d778=\"],\n  [\"88\u00d784=\", \"27\u00d742=\"],\n  [\"60\u00d730=\", \"36\u00d785=\"],\n  [\"91\u00d796=\", \"24\u00d761=\"],\n  [\"84\u00d717=\", \"69\u00d743=\"],\n  [\"65\u00d714=\", \"43\u00d799=\"],\n  [\"79\u00d781=\", \"28\u00d745=\"],\n  [\"61\u00d750=\", \"39\u00d767=\"],\n  [\"30\u00d734=\", \"77\u00d739=\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (const range of results.items) {\n    range.insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Replace the date header and all the two-digit multiplication problems.\n# Every source string in this document appears exactly once, so a\n# Find/Replace (wdReplaceAll) pass for each pair is safe and unambiguous.\n$d = $word.ActiveDocument\n\n$pairs = @(\n  @(\"2024-03-07 Thursday\", \"2024-03-08 Friday\"),\n  @(\"40\u00d748=\", \"13\u00d777=\"),\n  @(\"74\u00d798=\", \"76\u00d732=\"),\n  @(\"94\u00d788=\", \"76\u00d783=\"),\n  @(\"85\u00d711=\", \"32\u00d720=\"),\n  @(\"30\u00d762=\", \"71\u00d792=\"),\n  @(\"57\u00d785=\", \"83\u00d727=\"),\n  @(\"76\u00d748=\", \"28\u00d740=\"),\n  @(\"17\u00d776=\", \"87\u00d748=\"),\n  @(\"53\u00d748=\", \"28\u00d781=\"),\n  @(\"45\u00d789=\", \"29\u00d716=\"),\n  @(\"97\u00d748=\", \"92\u00d763=\"),\n  @(\"88\u00d715=\", \"66\u00d775=\"),\n  @(\"14\u00d741=\", \"84\u00d727=\"),\n  @(\"67\u00d722=\", \"53\u00d732=\"),\n  @(\"80\u00d793=\", \"18\u00d799=\"),\n  @(\"42\u00d734=\", \"38\u00d786=\"),\n  @(\"62\u00d777=\", \"61\u00d778=\"),\n  @(\"88\u00d784=\", \"27\u00d742=\"),\n  @(\"60\u00d730=\", \"36\u00d785=\"),\n  @(\"91\u00d796=\", \"24\u00d761=\"),\n  @(\"84\u00d717=\", \"69\u00d743=\"),\n  @(\"65\u00d714=\", \"43\u00d799=\"),\n  @(\"79\u00d781=\", \"28\u00d745=\"),\n  @(\"61\u00d750=\", \"39\u00d767=\"),\n  @(\"30\u00d734=\", \"77\u00d739=\")\n)\n\nforeach ($pair in $pairs) {\n  $old = $pair[0]\n  $new = $pair[1]\n  $find = $d.Content.Find\n  $find.Text = $old\n  $find.Replacement.Text = $new\n  $find.Execute($old, $false, $false, $false, $false, $false, $true, 1, $false, $new, 2)\n}\n"}
